# New weekly price record for Espinaca / Vega Modelo de Temuco.
# A new row is inserted at row 127 (pushing the existing rows 127-190 down
# to 128-191), and the freshly inserted row 127 is populated with the new
# observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 127, shifting rows 127:190 down
# to 128:191 (this also pushes the sheet's used range/dimension to R191).
$ws.Rows.Item(127).Insert()

# Populate the newly inserted row 127 with the new data point.
$ws.Cells.Item(127, 1).Value  = 10
$ws.Cells.Item(127, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(127, 3).Value  = "La Araucanía"
$ws.Cells.Item(127, 4).Value  = 44813
$ws.Cells.Item(127, 5).Value  = 9
$ws.Cells.Item(127, 6).Value  = 100112012
$ws.Cells.Item(127, 7).Value  = "Espinaca"
$ws.Cells.Item(127, 8).Value  = "Sin especificar"
$ws.Cells.Item(127, 9).Value  = "Primera"
$ws.Cells.Item(127, 10).Value = 10
$ws.Cells.Item(127, 11).Value = 13000
$ws.Cells.Item(127, 12).Value = 13000
$ws.Cells.Item(127, 13).Value = 13000
$ws.Cells.Item(127, 14).Value = "`$/docena de atados"
$ws.Cells.Item(127, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(127, 16).Value = 4333
$ws.Cells.Item(127, 17).Value = 3
$ws.Cells.Item(127, 18).Value = "Hortaliza"
